$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email addresses in column B (shared strings)
$ws.Range("B1").Value = "vaibhavzade802@gmail.com"
$ws.Range("B2").Value = "vaibhavzade159@gmail.com"

# Move the active selection to F5
$ws.Range("F5").Select()
